# Addresses import by sheet name #26
#
# The original workbook has a single sheet ("Sheet1") with a small
# Name/Age table. This change adds a second worksheet named "Test Sheet"
# (placed right after "Sheet1") that contains the same Name/Age table
# plus one additional row ("Adam", 20), and makes that new sheet the
# active tab.

$wb = $excel.ActiveWorkbook

# Existing worksheet with the sample data.
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet immediately after Sheet1 and name it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test Sheet"

# Copy the Name/Age table from Sheet1 into the new sheet.
$ws1.Range("A1:B3").Copy()
$ws2.Range("A1").PasteSpecial()

# Add the extra data row for the new sheet.
$ws2.Range("A4").Value = "Adam"
$ws2.Range("B4").Value = 20

# Update the selections to reflect the edited state, and make sure the
# new "Test Sheet" ends up as the active/selected sheet.
$ws1.Range("A1:B3").Select()
$ws2.Range("A5").Select()
